$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.169.59'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.839.50'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.48'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6313'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07514'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2932'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.27'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.93%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '1.830.35'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.000'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6705'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '82.81'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009343'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.023'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '29.181.60'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '2.082.15'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.61'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '224.21'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.150'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '160.09'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1405'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.523'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.506'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05960'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +14.86%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.166'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.071'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7513'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.859'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.143'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.682'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.232.82'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.774'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01793'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.574'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.25%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8962'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.005'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '102.44'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.979.70'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('B46').Value = 'XinFinNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.07949'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +18.43%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '66.03'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.42%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.00000000123'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4081'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.044'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.52%  '
